$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fix capitalization of existing codename (row 26: BIOFAB, University of Washington)
$ws.Range("D26").Value = "BUW"

# Append new BioFoundry rows
$newRows = @(
    @("Biofactorial", "Vancouver", "Canada", "BFC", 20, "CAN", 49.2, -123.6),
    @("Living Measurment Systems Foundry", "Gaithersburg, MD", "United States of America", "LMS", 21, "USA", 40, -77.2),
    @("Cyberbiofoundry", "Fort Collins, CO", "United States of America", "CBF", 22, "USA", 40, -105.5)
)

$startRow = 27
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
    $ws.Cells.Item($r, 8).Value = $row[7]
}

# Update active selection to K18 as in diff
$ws.Range("K18").Select()
